$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header labels for the new "team record" columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the existing header formatting (bold, border, centered) from AC1
# onto the new header cells so they match the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team record (Wins=75, Losses=86, Ties=0) for every data row (2-53).
for ($r = 2; $r -le 53; $r++) {
    $ws.Cells.Item($r, 30).Value = 75
    $ws.Cells.Item($r, 31).Value = 86
    $ws.Cells.Item($r, 32).Value = 0
}
